$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108; this shifts the existing rows
# 108..223 down to 109..224 (matches the new dimension A1:R224).
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with the new weekly record.
$ws.Cells.Item(108, 1).Value = 6
$ws.Cells.Item(108, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(108, 3).Value = "Metropolitana"
$ws.Cells.Item(108, 4).Value = 44539
$ws.Cells.Item(108, 5).Value = 13
$ws.Cells.Item(108, 6).Value = 100112026
$ws.Cells.Item(108, 7).Value = "Haba"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 320
$ws.Cells.Item(108, 11).Value = 9000
$ws.Cells.Item(108, 12).Value = 10000
$ws.Cells.Item(108, 13).Value = 9375
$ws.Cells.Item(108, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(108, 15).Value = "Región del Maule"
$ws.Cells.Item(108, 16).Value = 375
$ws.Cells.Item(108, 17).Value = 25
$ws.Cells.Item(108, 18).Value = "Hortaliza"

# Restore the date-style (numFmt) on D108 that Insert() should have
# carried over automatically; set explicitly in case it didn't.
$ws.Cells.Item(108, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat
